# Update column G ("K") with newly computed strikeout (K) values,
# replacing the previous placeholder "Strike#" derived figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 3
    5  = 3
    6  = 1
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 3
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 4
    26 = 2
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
